$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 inherits the original E1 formatting (full red box + quote-prefixed blank).
$ws.Range("E1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# C1 and E1 inherit the existing "left+right+top" red box (no bottom) border,
# currently used by E3, before we change E3 itself.
$ws.Range("E3").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)

# Fill in the new legend values in row 1.
$ws.Range("C1").Value = " "
$ws.Range("D1").Value = "'"
$ws.Range("E1").Value = "'　"

# B1 loses its right edge, keeping only the top border.
$ws.Range("B1").Borders.Item(10).LineStyle = -4142

# Move the active selection from H6 to G6.
$ws.Range("G6").Select()
